$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68, shifting rows 68:160 down to 69:161
$ws.Rows("68:68").Insert()

# Populate the new row 68 with data
$ws.Range("A68").Value = 11
$ws.Range("B68").Value = "Vega Monumental Concepción"
$ws.Range("C68").Value = "Bíobío"
$ws.Range("D68").Value = 44848
$ws.Range("E68").Value = 8
$ws.Range("F68").Value = 100112043
$ws.Range("G68").Value = "Pepino ensalada"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 250
$ws.Range("K68").Value = 24000
$ws.Range("L68").Value = 25000
$ws.Range("M68").Value = 24400
$ws.Range("N68").Value = "$/caja 60 unidades"
$ws.Range("O68").Value = "Región de Arica y Parinacota"
$ws.Range("P68").Value = 407
$ws.Range("Q68").Value = 60
$ws.Range("R68").Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Range("D68").NumberFormat = $ws.Range("D69").NumberFormat
